$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 214, shifting existing rows 214:292 down to 215:293
$ws.Rows.Item(214).Insert()

# Populate the newly inserted row 214 with data
$ws.Cells.Item(214, 1).Value = 11
$ws.Cells.Item(214, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(214, 3).Value = "Bíobío"
$ws.Cells.Item(214, 4).Value = 44875
$ws.Cells.Item(214, 5).Value = 8
$ws.Cells.Item(214, 6).Value = 100114013
$ws.Cells.Item(214, 7).Value = "Zanahoria"
$ws.Cells.Item(214, 8).Value = "Sin especificar"
$ws.Cells.Item(214, 9).Value = "Primera"
$ws.Cells.Item(214, 10).Value = 450
$ws.Cells.Item(214, 11).Value = 13000
$ws.Cells.Item(214, 12).Value = 14000
$ws.Cells.Item(214, 13).Value = 13556
$ws.Cells.Item(214, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(214, 15).Value = "Región Metropolitana"
$ws.Cells.Item(214, 16).Value = 678
$ws.Cells.Item(214, 17).Value = 20
$ws.Cells.Item(214, 18).Value = "Hortaliza"
